$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bug entries added to the "Gameplay" sheet, column A only (rows 10-17)
$newBugs = @(
    "Knight Skill 4 Spear too high+",
    "Knight attack speed too fast",
    "Infantry Skill 1 Hook exception",
    "Fall through not working with controller ",
    "Fall through causes stick in ground controller",
    "Champion Selection: when 3 players, player 1 and 3 target same preview platform",
    "Wrong team shown as winning",
    "Running continues after game ended"
)

$row = 10
foreach ($bug in $newBugs) {
    $ws.Cells.Item($row, 1).Value = $bug
    $row = $row + 1
}

# Select the last filled cell, matching the saved selection state
$ws.Range("A17").Select()
